# OpenEMRData.xlsx edit: add "bala/Danish" row to invalidCredentialTest and
# add a new "addPatientTest" worksheet with patient-add test data.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- invalidCredentialTest: append row 5 ---
$ws1.Range("A5").Value = "bala"
$ws1.Range("B5").Value = "bala123"
$ws1.Range("C5").Value = "Danish"
$ws1.Range("D5").Value = "1Invalid username or password2"

# --- add new worksheet "addPatientTest" right after invalidCredentialTest ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "addPatientTest"

# Header row
$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "Password"
$ws2.Range("C1").Value = "Language"
$ws2.Range("D1").Value = "FirstName"
$ws2.Range("E1").Value = "LastName"
$ws2.Range("F1").Value = "Dob"
$ws2.Range("G1").Value = "Gender"
$ws2.Range("H1").Value = "ExpectedAlertText"
$ws2.Range("I1").Value = "ExpectedValue"

# The "Dob" column is stored as text (numFmtId 49 / "@") so the date-looking
# string isn't reinterpreted as a serial date value.
$ws2.Range("F1:F2").NumberFormat = "@"

# Data row
$ws2.Range("A2").Value = "admin"
$ws2.Range("B2").Value = "pass"
$ws2.Range("C2").Value = "English (Indian)"
$ws2.Range("D2").Value = "Bala"
$ws2.Range("E2").Value = "dina"
$ws2.Range("F2").Value = "2021-07-20"
$ws2.Range("G2").Value = "Male"
$ws2.Range("H2").Value = "Tobacco"
$ws2.Range("I2").Value = "Medical Record Dashboard - Bala Dina"

# Column widths (bestFit-style, matching the authored workbook's column sizing)
$ws2.Columns.Item(1).ColumnWidth = 9.166666666666668
$ws2.Columns.Item(2).ColumnWidth = 8.666666666666668
$ws2.Columns.Item(3).ColumnWidth = 14.0
$ws2.Columns.Item(4).ColumnWidth = 9.333333333333332
$ws2.Columns.Item(5).ColumnWidth = 8.833333333333332
$ws2.Columns.Item(6).ColumnWidth = 19.166666666666664
$ws2.Columns.Item(7).ColumnWidth = 6.666666666666666
$ws2.Columns.Item(8).ColumnWidth = 16.833333333333336
$ws2.Columns.Item(9).ColumnWidth = 34.16666666666667

# Print orientation
$ws2.PageSetup.Orientation = 1

# Selections / active sheet, matching the authored file
$ws1.Range("D5").Select()
$ws2.Range("B5").Select()
$ws2.Activate()
